# Normalize the "Recorded By" (column G) cell values on the
# "Session Analysis Results" sheet so that each comma-separated list of
# names/emails is sorted: entries equal to "System" (case-insensitively)
# come first (ties broken by ordinal/ASCII order, so "System" sorts
# before "system"), followed by the remaining entries in ordinal order.

function Compare-Names($x, $y) {
    $xIsSystem = ($x.ToLower() -eq "system")
    $yIsSystem = ($y.ToLower() -eq "system")
    if ($xIsSystem -and -not $yIsSystem) { return -1 }
    if ($yIsSystem -and -not $xIsSystem) { return 1 }
    return $x.CompareTo($y)
}

function Sort-Names($arr) {
    # simple, dependency-free insertion sort using Compare-Names
    $n = $arr.Count
    for ($i = 1; $i -lt $n; $i++) {
        $key = $arr[$i]
        $j = $i - 1
        while ($j -ge 0 -and (Compare-Names $arr[$j] $key) -gt 0) {
            $arr[$j + 1] = $arr[$j]
            $j = $j - 1
        }
        $arr[$j + 1] = $key
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data occupies rows 2..157 (row 1 is the header row).
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text
    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text -split ',\s*'
    if ($parts.Count -le 1) { continue }

    $sortedParts = Sort-Names $parts
    $newText = [string]::Join(", ", $sortedParts)

    if (-not $newText.Equals($text)) {
        $cell.Value = $newText
    }
}
